$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5906755328178406
$ws.Range("B1").Value = 1.473466873168945
$ws.Range("C1").Value = 5.993953227996826
$ws.Range("D1").Value = 1.93631386756897
$ws.Range("E1").Value = 1.551906108856201
